$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 with same style as E1 (bold, centered, bordered)
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row (no special style, like column B-E data cells)
$ws.Range("F2").Value = "2021-10-05 10:51:31.665286"
$ws.Range("F3").Value = "2021-10-05 10:51:31.665299"
$ws.Range("F4").Value = "2021-10-05 10:51:31.665303"
$ws.Range("F5").Value = "2021-10-05 10:51:31.665306"
$ws.Range("F6").Value = "2021-10-05 10:51:31.665310"
$ws.Range("F7").Value = "2021-10-05 10:51:31.665313"
$ws.Range("F8").Value = "2021-10-05 10:51:31.665316"
$ws.Range("F9").Value = "2021-10-05 10:51:31.665319"
$ws.Range("F10").Value = "2021-10-05 10:51:31.665322"
$ws.Range("F11").Value = "2021-10-05 10:51:31.665325"
$ws.Range("F12").Value = "2021-10-05 10:51:31.665328"
$ws.Range("F13").Value = "2021-10-05 10:51:31.665331"
$ws.Range("F14").Value = "2021-10-05 10:51:31.665334"
$ws.Range("F15").Value = "2021-10-05 10:51:31.665337"
$ws.Range("F16").Value = "2021-10-05 10:51:31.665340"
$ws.Range("F17").Value = "2021-10-05 10:51:31.665343"
$ws.Range("F18").Value = "2021-10-05 10:51:31.665346"
$ws.Range("F19").Value = "2021-10-05 10:51:31.665350"
$ws.Range("F20").Value = "2021-10-05 10:51:31.665353"
$ws.Range("F21").Value = "2021-10-05 10:51:31.665356"

Write-Host "Added time_taken column (F1:F21)"
